# Update NATMI LR-pair values per author revision ("Natmi following Dr Hou advice")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.912301666666667
$ws.Range("H2").Value = 5.736905
$ws.Range("I2").Value = 0.1465770754282357
$ws.Range("J2").Value = 0.1717732502851657
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.690981
$ws.Range("N2").Value = 2.072943
$ws.Range("O2").Value = 0.8010798064371025
$ws.Range("P2").Value = 0.8579689210526861
$ws.Range("Q2").Value = 1.321364117935
$ws.Range("R2").Value = 11.892277061415
$ws.Range("S2").Value = 0.1174199352121676
$ws.Range("T2").Value = 0.1473761102128766
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.912301666666667
$ws.Range("H3").Value = 5.736905
$ws.Range("I3").Value = 0.1465770754282357
$ws.Range("J3").Value = 0.1717732502851657
$ws.Range("K3").Value = 2
$ws.Range("M3").Value = 0.171581
$ws.Range("N3").Value = 0.343162
$ws.Range("O3").Value = 0.1989201935628975
$ws.Range("P3").Value = 0.142031078947314
$ws.Range("Q3").Value = 0.3281146322683334
$ws.Range("R3").Value = 1.96868779361
$ws.Range("S3").Value = 0.02915714021606807
$ws.Range("T3").Value = 0.02439714007228909
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.724098666666666
$ws.Range("H4").Value = 14.172296
$ws.Range("I4").Value = 0.3621000695990751
$ws.Range("J4").Value = 0.4243440231141098
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.690981
$ws.Range("N4").Value = 2.072943
$ws.Range("O4").Value = 0.8010798064371025
$ws.Range("P4").Value = 0.8579689210526861
$ws.Range("Q4").Value = 3.264262420791999
$ws.Range("R4").Value = 29.378361787128
$ws.Range("S4").Value = 0.2900710536652885
$ws.Range("T4").Value = 0.3640739836663689
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.724098666666666
$ws.Range("H5").Value = 14.172296
$ws.Range("I5").Value = 0.3621000695990751
$ws.Range("J5").Value = 0.4243440231141098
$ws.Range("K5").Value = 2
$ws.Range("M5").Value = 0.171581
$ws.Range("N5").Value = 0.343162
$ws.Range("O5").Value = 0.1989201935628975
$ws.Range("P5").Value = 0.142031078947314
$ws.Range("Q5").Value = 0.8105655733253333
$ws.Range("R5").Value = 4.863393439952
$ws.Range("S5").Value = 0.07202901593378669
$ws.Range("T5").Value = 0.06027003944774097
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.3397296666666667
$ws.Range("H6").Value = 1.019189
$ws.Range("I6").Value = 0.02604012841917865
$ws.Range("J6").Value = 0.03051635109608539
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.690981
$ws.Range("N6").Value = 2.072943
$ws.Range("O6").Value = 0.8010798064371025
$ws.Range("P6").Value = 0.8579689210526861
$ws.Range("Q6").Value = 0.234746744803
$ws.Range("R6").Value = 2.112720703227
$ws.Range("S6").Value = 0.02086022103363292
$ws.Range("T6").Value = 0.02618208082437334
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.3397296666666667
$ws.Range("H7").Value = 1.019189
$ws.Range("I7").Value = 0.02604012841917865
$ws.Range("J7").Value = 0.03051635109608539
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 0.171581
$ws.Range("N7").Value = 0.343162
$ws.Range("O7").Value = 0.1989201935628975
$ws.Range("P7").Value = 0.142031078947314
$ws.Range("Q7").Value = 0.05829115593633334
$ws.Range("R7").Value = 0.3497469356180001
$ws.Range("S7").Value = 0.005179907385545726
$ws.Range("T7").Value = 0.004334270271712056
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.3292183333333333
$ws.Range("H8").Value = 0.987655
$ws.Range("I8").Value = 0.02523443937664543
$ws.Range("J8").Value = 0.02957216643998729
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.690981
$ws.Range("N8").Value = 2.072943
$ws.Range("O8").Value = 0.8010798064371025
$ws.Range("P8").Value = 0.8579689210526861
$ws.Range("Q8").Value = 0.227483613185
$ws.Range("R8").Value = 2.047352518665
$ws.Range("S8").Value = 0.02021479981139192
$ws.Range("T8").Value = 0.02537199973370635
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.3292183333333333
$ws.Range("H9").Value = 0.987655
$ws.Range("I9").Value = 0.02523443937664543
$ws.Range("J9").Value = 0.02957216643998729
$ws.Range("K9").Value = 2
$ws.Range("M9").Value = 0.171581
$ws.Range("N9").Value = 0.343162
$ws.Range("O9").Value = 0.1989201935628975
$ws.Range("P9").Value = 0.142031078947314
$ws.Range("Q9").Value = 0.05648761085166667
$ws.Range("R9").Value = 0.33892566511
$ws.Range("S9").Value = 0.005019639565253513
$ws.Range("T9").Value = 0.004200166706280944
$ws.Range("E10").Value = 2
$ws.Range("G10").Value = 5.7410415
$ws.Range("H10").Value = 11.482083
$ws.Range("I10").Value = 0.4400482871768651
$ws.Range("J10").Value = 0.3437942090646517
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.690981
$ws.Range("N10").Value = 2.072943
$ws.Range("O10").Value = 0.8010798064371025
$ws.Range("P10").Value = 0.8579689210526861
$ws.Range("Q10").Value = 3.9669505967115
$ws.Range("R10").Value = 23.801703580269
$ws.Range("S10").Value = 0.3525137967146216
$ws.Range("T10").Value = 0.2949647466153608
$ws.Range("E11").Value = 2
$ws.Range("G11").Value = 5.7410415
$ws.Range("H11").Value = 11.482083
$ws.Range("I11").Value = 0.4400482871768651
$ws.Range("J11").Value = 0.3437942090646517
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 0.171581
$ws.Range("N11").Value = 0.343162
$ws.Range("O11").Value = 0.1989201935628975
$ws.Range("P11").Value = 0.142031078947314
$ws.Range("Q11").Value = 0.9850536416115
$ws.Range("R11").Value = 3.940214566446
$ws.Range("S11").Value = 0.08753449046224353
$ws.Range("T11").Value = 0.04882946244929092
